$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column header (H1): "actually measured" ---
# Leading apostrophe forces text storage (shared string) instead of Excel's
# automatic numeric coercion; ClearFormats() strips the resulting
# quote-prefix cell style again so no stray formatting is left behind.
$ws.Range("H1").Value = "'actually measured"
$ws.Range("H1").ClearFormats()

# --- Row 2 (XL_SPT): replace design numbers with the actually-measured
#     ones, and keep the old design numbers as a reference block in H2:K2 ---
$ws.Range("B2").Value = "'4.65"
$ws.Range("C2").Value = "'9.5"
$ws.Range("D2").Value = "'0.73"
$ws.Range("E2").Value = "'2.35"
$ws.Range("B2:E2").ClearFormats()

$ws.Range("H2").Value = 4.7

$ws.Range("I2").Value = "'9.55"
$ws.Range("I2").ClearFormats()

$ws.Range("J2").Value = "'0.78"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").Value = 2.4

# --- View state: selection moved from E7 to G8 ---
$ws.Range("G8").Select()
